$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Split column A off from the merged col(1:2) width definition so it gets its own
#    <col min="1" max="1"/> entry (column B already has its own separate override).
$ws.Columns.Item(1).Hidden = $false

# 2. Update the "Objetivos:" row (row 10) with the new objectives paragraph.
$ws.Range("B10:C10").Value = "Apresentar os princípios fundamentais envolvidos nas operações relacionadas a sistemas particulados, de forma a permitir a análise de desempenho dos equipamentos que lidam com estes sistemas."

# 3. Insert a new blank row at position 13 (shifts old rows 13-23 down to 14-24,
#    preserving their row heights/formatting).
$ws.Rows("13:13").Insert()

# 4. New row 13: professor name, moved here from the old "Metodo" row mistake.
$ws.Range("B13:C13").Value = "787307 - Luis Fernando Figueiredo Faria"

# 5. Row 14 ("Programa resumido:") gets the new short-syllabus text.
$ws.Range("B14:C14").Value = "Fundamentos e caracterização de partículas e sistemas particulados. Dinâmica da interação sólido-fluido. Aplicações em sistemas diluídos: elutriação, câmara de poeira, ciclones, centrífugas e hidrociclones. Aplicações em sistemas concentrados: escoamento monofásico em meios porosos, filtração sólido-líquido, sedimentação, fluidização, transporte pneumático e hidráulico de partículas."

# 6. Row 16 ("Programa:") gets the new detailed syllabus text.
$ws.Range("B16:C16").Value = "1. Caracterização de partículas e sistemas particulados: noções de amostragem; diâmetros de esferas equivalentes e  diâmetros estatísticos; esfericidade; análise granulométrica, frequência simples e acumuladas; modelos de distribuição de tamanhos.   `n2. Interação sólido-fluido: Dinâmica e análise dimensional do sistema partícula - fluido infinito:velocidade terminal; lei de Stokes; correlação entre coeficiente de arraste e número de Reynolds para esferas; efeito da forma das partículas; efeito de paredes; efeito de população; efeito de deslizamento. `n3. Aplicações em sistemas diluídos: separação sólido-sólido por elutriação; partículas equitombantes e razão de sedimentação; eficiências globais e individuais de coleta; diâmetro de corte; separação sólido-gás com câmaras de poeira e ciclones; separação sólido-líquido com centrífugas e hidrociclones.`n4. Aplicações em  sistemas concentrados: escoamento monofásico em meios porosos; separação sólido-líquido por filtração em superfície; auxiliares de filtração; estudo detalhado dos filtros prensa e de tambor rotativo; separação sólido-líquido por sedimentação; leitos fluidizados a gás e a líquido; curva característica e histerese de fluidização; previsão das velocidades mínima e máxima de fluidização; transporte pneumático de partículas; velocidade de deslizamento; transporte hidráulico de partículas; velocidade de salto."

# 7. Row 19 ("Metodo:") gets the evaluation method text.
$ws.Range("B19:C19").Value = "Participação em sala de aula, preparação e apresentação de trabalhos e provas escritas."

# 8. Row 20 ("Criterio:") gets the final-grade formula text.
$ws.Range("B20:C20").Value = "Média Final = (Prova1 + Prova2 + Nota de Trabalho) /3`nMédia final mínima de aprovação = 5,0"

# 9. Row 21 ("Norma de recuperacao:") gets the make-up exam formula text.
$ws.Range("B21:C21").Value = "(Prova escrita + Média Final)/2         Nota Final mínima para aprovação= 5,0"

# 10. Row 22 ("Bibliografia:") gets the new bibliography list.
$ws.Range("B22:C22").Value = "1. PERRY, R.H.; GREEN, D.W.; MALONEY, J.O. (Eds.). Perrys Chemical Engineers Handbook. New York : McGraw-Hill, 1997.`n2. MASSARANI, G. Fluidodinâmica em Sistemas Particulados. 2. ed. RJ: E-Papers, 2002.`n3. SVAROVSKY, L. Solid-Liquid Separation. 3. ed. LondonBoston: Butterworths, 1990.`n4. RUSHTON, A.; WARD, A.S.; HOLDICH, R.G. Solid-Liquid Filtration and Separation Technology. Weinheim:  VCH, 1996.`n5. COULSON, J.M.; RICHARDSON, J.F. Chemical Engineering. 5th. ed. Londres: Pergamon Press,1996. Vol. 2.`n6. ALLEN, T. Particle Size Measurement. 5th. ed. Londres: Chapman & Hall, 1997. Vol 1 e 2."

